$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "zvihcxioc" -> split into two runs "Z" / "vihcxioc" (first letter
#    capitalized).  We briefly toggle Bold on/off on the single leading
#    character so the engine keeps the run boundary it creates while editing
#    that character, instead of silently re-coalescing it back into the
#    neighbouring run of identical formatting once we are done.
# ---------------------------------------------------------------------------
$firstChar = $d.Range(0, 1)
$firstChar.Text = "Z"

$firstCharAgain = $d.Range(0, 1)
$firstCharAgain.Bold = $true
$firstCharAgain.Bold = $false

# ---------------------------------------------------------------------------
# 2) Add a new line/paragraph containing "ftfhjghfgdtu" right where the
#    existing _GoBack bookmark sits, so the bookmark ends up wrapping the new
#    paragraph exactly like in the target document.
# ---------------------------------------------------------------------------
$bookmark = $d.Bookmarks.Item("_GoBack")
$breakPoint = $d.Range($bookmark.Start, $bookmark.Start)
$breakPoint.InsertBefore("`r")

$bookmarkAfterBreak = $d.Bookmarks.Item("_GoBack")
$newLineStart = $d.Range($bookmarkAfterBreak.Start, $bookmarkAfterBreak.Start)
$newLineStart.InsertBefore("ftfhjghfgdtu")

Write-Host "Paragraphs:" $d.Paragraphs.Count
foreach ($p in $d.Paragraphs) {
    Write-Host ("  [{0}]" -f $p.Range.Text)
}
